# [Feat 2269][Feat 3238] Added custom coercer for ImportMode enum cells.
# This script mutates the TEST_CASES sheet of the import template:
#  - removes the REPLACE/R demo rows (old rows 4 & 5)
#  - renumbers the path/ref/name sample values so they stay sequential
#  - populates the previously-empty TC_DESCRIPTION, TC_PRE_REQUISITE,
#    TC_CREATED_ON and TC_CREATED_BY sample columns
#  - refreshes the sheet view so the newly-filled columns are visible

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the two "REPLACE"/"R" sample rows; everything below shifts up.
$ws.Rows("4:5").Delete() | Out-Null

# Renumber the TC_PATH / TC_REFERENCE / TC_NAME / TC_NUM sample data so the
# sequence (1..8) is contiguous again after the row removal.
for ($i = 3; $i -le 8; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 4).Value2 = "path/row$i"
    $ws.Cells.Item($row, 5).Value2 = 10 + $i
    $ws.Cells.Item($row, 7).Value2 = "ref$i"
    $ws.Cells.Item($row, 8).Value2 = "name$i"
}

# Fill in the new TC_DESCRIPTION (N), TC_PRE_REQUISITE (O), TC_CREATED_ON (S)
# and TC_CREATED_BY (T) sample values for all 8 data rows.
$baseDate = Get-Date -Year 2003 -Month 2 -Day 1
for ($i = 1; $i -le 8; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 14).Value2 = "desc$i"
    $ws.Cells.Item($row, 15).Value2 = "pre$i"
    $ws.Cells.Item($row, 19).Value2 = $baseDate.AddDays($i - 1)
    $ws.Cells.Item($row, 19).NumberFormat = "m/d/yyyy"
    $ws.Cells.Item($row, 20).Value2 = "creator$i"
}

# Update the sheet view to reflect the user scrolling to inspect the new
# TC_CREATED_BY column.
$ws.Range("T2:T9").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 14
